# Rename "Disturbia.wav" to "Disturbia.mp3" and make it the active sheet
# with D20 selected (matches the author's edits in the commit
# "updated files for disturbia.mp3").

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Disturbia.wav")
$ws.Name = "Disturbia.mp3"

$ws.Activate()
$ws.Range("D20").Select()
